# MSAA Enable and update rubric
# Fill in the "Milestone Complete(X)" flags (column F) and the student/git
# rows that were left blank, plus the two citation rows at the bottom of
# the sheet (A95/A96). Excel recalculates every dependent formula
# (G/H/I/J/K/L columns) automatically once the inputs change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Milestone Complete(X) marks for the GEOMETRY / CAMERA & VIEW rows
$ws.Range("F4").Value  = "X"
$ws.Range("F5").Value  = "X"
$ws.Range("F18").Value = "X"
$ws.Range("F23").Value = "X"
$ws.Range("F30").Value = "X"
$ws.Range("F31").Value = "X"
$ws.Range("F32").Value = "X"
$ws.Range("F33").Value = "X"

# MSAA rows (54/55) now record which milestone the feature shipped in,
# plus the completion flag
$ws.Range("E54").Value = "I"
$ws.Range("F54").Value = "X"
$ws.Range("E55").Value = "I"
$ws.Range("F55").Value = "X"

# Effective use of GIT / cleaned up graphics objects, milestone I complete
$ws.Range("C90").Value = "X"
$ws.Range("C91").Value = "X"

# Project source citation
$ws.Range("A95").Value = "http://www.braynzarsoft.net/"
$ws.Range("A96").Value = "3D Game Programming with DirectX 11"

# Leave the selection where the author finished editing
$ws.Activate()
$ws.Range("F32").Select()
